{"js": "// Office.js (Word JavaScript API) edit script.\n// The document consists of one title paragraph (\"<date> <weekday>\")\n// followed by a single table of 20 rows x 5 columns (100 cells), where\n// every cell holds exactly one paragraph with exactly one run of text\n// (an arithmetic expression like \"47+48=95\"). `context.document.body\n// .paragraphs` enumerates ALL paragraphs in document order, including\n// those nested inside table cells, so the title paragraph is index 0\n// and the table cells follow in row-major order (row0 col0..col4,\n// row1 col0..col4, ...) \u2014 this lines up exactly with the unified diff,\n// which replaces each <w:t> run's text (still index 0 = title, 1..100\n// = cells) with a new value while preserving every run's formatting\n// (rFonts/sz) untouched.\nconst REPLACEMENTS = [\n  [\"2025-01-25 Saturday\", \"2025-01-26 Sunday\"],\n  [\"47+48=95\", \"57-56=1\"],\n  [\"90-67=23\", \"51+34=85\"],\n  [\"32+37=69\", \"20+70=90\"],\n  [\"20+9=29\", \"25+49=74\"],\n  [\"74-50=24\", \"58+18=76\"],\n  [\"86-5=81\", \"90-15=75\"],\n  [\"34+49=83\", \"29+9=38\"],\n  [\"38+39=77\", \"29+30=59\"],\n  [\"29+19=48\", \"10+88=98\"],\n  [\"82-77=5\", \"81-77=4\"],\n  [\"26+59=85\", \"97-21=76\"],\n  [\"2+96=98\", \"40+33=73\"],\n  [\"80-3=77\", \"71+21=92\"],\n  [\"51+3=54\", \"72-72=0\"],\n  [\"36+33=69\", \"41+52=93\"],\n  [\"87-41=46\", \"61+11=72\"],\n  [\"6+32=38\", \"28+21=49\"],\n  [\"29-8=21\", \"39+22=61\"],\n  [\"66-30=36\", \"80+14=94\"],\n  [\"28+49=77\", \"60-8=52\"],\n  [\"77-53=24\", \"91-37=54\"],\n  [\"90-32=58\", \"7+7=14\"],\n  [\"70+26=96\", \"65+15=80\"],\n  [\"5+2=7\", \"54+12=66\"],\n  [\"23+6=29\", \"77+15=92\"],\n  [\"43-11=32\", \"7+11=18\"],\n  [\"66+30=96\", \"40+15=55\"],\n  [\"88-55=33\", \"10+89=99\"],\n  [\"31-22=9\", \"24+13=37\"],\n  [\"99-79=20\", \"7-5=2\"],\n  [\"8+56=64\", \"38-19=19\"],\n  [\"75-65=10\", \"99-78=21\"],\n  [\"85+1=86\", \"23+28=51\"],\n  [\"81-76=5\", \"91-20=71\"],\n  [\"25+16=41\", \"47+11=58\"],\n  [\"32+47=79\", \"6+60=66\"],\n  [\"83+10=93\", \"85-55=30\"],\n  [\"76-24=52\", \"30-12=18\"],\n  [\"51-47=4\", \"39+42=81\"],\n  [\"56+41=97\", \"12-5=7\"],\n  [\"67-13=54\", \"43+7=50\"],\n  [\"1+82=83\", \"69-29=40\"],\n  [\"73-19=54\", \"49-22=27\"],\n  [\"9+24=33\", \"43-22=21\"],\n  [\"7+12=19\", \"5+9=14\"],\n  [\"8+63=71\", \"9+32=41\"],\n  [\"65-59=6\", \"65-50=15\"],\n  [\"57+35=92\", \"13+55=68\"],\n  [\"24+12=36\", \"88-51=37\"],\n  [\"60-18=42\", \"55-16=39\"],\n  [\"9+14=23\", \"81-3=78\"],\n  [\"28+43=71\", \"89-77=12\"],\n  [\"17-4=13\", \"66+27=93\"],\n  [\"21+70=91\", \"0+76=76\"],\n  [\"17-1=16\", \"50+24=74\"],\n  [\"82-44=38\", \"49-29=20\"],\n  [\"39+13=52\", \"46-9=37\"],\n  [\"23+26=49\", \"84-11=73\"],\n  [\"10+42=52\", \"32+22=54\"],\n  [\"64-41=23\", \"93-40=53\"],\n  [\"81+7=88\", \"76-20=56\"],\n  [\"13+5=18\", \"66-9=57\"],\n  [\"59-1=58\", \"98-76=22\"],\n  [\"65-57=8\", \"62-54=8\"],\n  [\"27-21=6\", \"8-5=3\"],\n  [\"78-30=48\", \"47+25=72\"],\n  [\"92-44=48\", \"52+42=94\"],\n  [\"70-35=35\", \"68-51=17\"],\n  [\"66-23=43\", \"60+29=89\"],\n  [\"7+57=64\", \"19-15=4\"],\n  [\"88-83=5\", \"1+56=57\"],\n  [\"91-38=53\", \"84-68=16\"],\n  [\"58+39=97\", \"43-2=41\"],\n  [\"65-21=44\", \"96-77=19\"],\n  [\"65+6=71\", \"81-37=44\"],\n  [\"80-39=41\", \"19+50=69\"],\n  [\"16+20=36\", \"1+37=38\"],\n  [\"3+86=89\", \"36+20=56\"],\n  [\"23+49=72\", \"19+53=72\"],\n  [\"15+2=17\", \"80-15=65\"],\n  [\"13+25=38\", \"60-53=7\"],\n  [\"86-64=22\", \"13+78=91\"],\n  [\"30+52=82\", \"68-22=46\"],\n  [\"88-71=17\", \"83+9=92\"],\n  [\"31-28=3\", \"64-53=11\"],\n  [\"43-43=0\", \"35+55=90\"],\n  [\"45+27=72\", \"37+61=98\"],\n  [\"90-57=33\", \"67+17=84\"],\n  [\"56-12=44\", \"68-16=52\"],\n  [\"89-51=38\", \"66+33=99\"],\n  [\"36+48=84\", \"21+41=62\"],\n  [\"94-63=31\", \"46+50=96\"],\n  [\"71-40=31\", \"29+59=88\"],\n  [\"6+77=83\", \"38+61=99\"],\n  [\"9+54=63\", \"99-49=50\"],\n  [\"28-21=7\", \"37-8=29\"],\n  [\"75-53=22\", \"29-11=18\"],\n  [\"64+9=73\", \"11+42=53\"],\n  [\"5+70=75\", \"59+7=66\"],\n  [\"41-33=8\", \"58+10=68\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    `Expected ${REPLACEMENTS.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < REPLACEMENTS.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = paragraphs.items[i];\n  // Defensive check: only rewrite the run if it still holds the text we\n  // expect from the \"before\" snapshot; otherwise leave it untouched so an\n  // already-edited / reordered document doesn't get silently mangled.\n  if (para.text === oldText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) edit script.\n# The document has one title paragraph (a date/weekday line) followed\n# by a single 20-row x 5-column table whose 100 cells each hold exactly\n# one paragraph with one run of text (an arithmetic expression, e.g.\n# \"47+48=95\"). We rewrite the title paragraph and every cell in row-major\n# order (row0 col0..col4, row1 col0..col4, ...), matching the unified diff,\n# via Range.Text assignment so each run's existing formatting\n# (rFonts/sz) is preserved and only the <w:t> contents change.\n\n$DateOld = '2025-01-25 Saturday'\n$DateNew = '2025-01-26 Sunday'\n\n# Row-major (row0col0..row0col4, row1col0.., ...) old/new text pairs\n# for the 20x5 answer table, taken from the unified diff in document order.\n$CellPairs = @(\n  @('47+48=95', '57-56=1'),\n  @('90-67=23', '51+34=85'),\n  @('32+37=69', '20+70=90'),\n  @('20+9=29', '25+49=74'),\n  @('74-50=24', '58+18=76'),\n  @('86-5=81', '90-15=75'),\n  @('34+49=83', '29+9=38'),\n  @('38+39=77', '29+30=59'),\n  @('29+19=48', '10+88=98'),\n  @('82-77=5', '81-77=4'),\n  @('26+59=85', '97-21=76'),\n  @('2+96=98', '40+33=73'),\n  @('80-3=77', '71+21=92'),\n  @('51+3=54', '72-72=0'),\n  @('36+33=69', '41+52=93'),\n  @('87-41=46', '61+11=72'),\n  @('6+32=38', '28+21=49'),\n  @('29-8=21', '39+22=61'),\n  @('66-30=36', '80+14=94'),\n  @('28+49=77', '60-8=52'),\n  @('77-53=24', '91-37=54'),\n  @('90-32=58', '7+7=14'),\n  @('70+26=96', '65+15=80'),\n  @('5+2=7', '54+12=66'),\n  @('23+6=29', '77+15=92'),\n  @('43-11=32', '7+11=18'),\n  @('66+30=96', '40+15=55'),\n  @('88-55=33', '10+89=99'),\n  @('31-22=9', '24+13=37'),\n  @('99-79=20', '7-5=2'),\n  @('8+56=64', '38-19=19'),\n  @('75-65=10', '99-78=21'),\n  @('85+1=86', '23+28=51'),\n  @('81-76=5', '91-20=71'),\n  @('25+16=41', '47+11=58'),\n  @('32+47=79', '6+60=66'),\n  @('83+10=93', '85-55=30'),\n  @('76-24=52', '30-12=18'),\n  @('51-47=4', '39+42=81'),\n  @('56+41=97', '12-5=7'),\n  @('67-13=54', '43+7=50'),\n  @('1+82=83', '69-29=40'),\n  @('73-19=54', '49-22=27'),\n  @('9+24=33', '43-22=21'),\n  @('7+12=19', '5+9=14'),\n  @('8+63=71', '9+32=41'),\n  @('65-59=6', '65-50=15'),\n  @('57+35=92', '13+55=68'),\n  @('24+12=36', '88-51=37'),\n  @('60-18=42', '55-16=39'),\n  @('9+14=23', '81-3=78'),\n  @('28+43=71', '89-77=12'),\n  @('17-4=13', '66+27=93'),\n  @('21+70=91', '0+76=76'),\n  @('17-1=16', '50+24=74'),\n  @('82-44=38', '49-29=20'),\n  @('39+13=52', '46-9=37'),\n  @('23+26=49', '84-11=73'),\n  @('10+42=52', '32+22=54'),\n  @('64-41=23', '93-40=53'),\n  @('81+7=88', '76-20=56'),\n  @('13+5=18', '66-9=57'),\n  @('59-1=58', '98-76=22'),\n  @('65-57=8', '62-54=8'),\n  @('27-21=6', '8-5=3'),\n  @('78-30=48', '47+25=72'),\n  @('92-44=48', '52+42=94'),\n  @('70-35=35', '68-51=17'),\n  @('66-23=43', '60+29=89'),\n  @('7+57=64', '19-15=4'),\n  @('88-83=5', '1+56=57'),\n  @('91-38=53', '84-68=16'),\n  @('58+39=97', '43-2=41'),\n  @('65-21=44', '96-77=19'),\n  @('65+6=71', '81-37=44'),\n  @('80-39=41', '19+50=69'),\n  @('16+20=36', '1+37=38'),\n  @('3+86=89', '36+20=56'),\n  @('23+49=72', '19+53=72'),\n  @('15+2=17', '80-15=65'),\n  @('13+25=38', '60-53=7'),\n  @('86-64=22', '13+78=91'),\n  @('30+52=82', '68-22=46'),\n  @('88-71=17', '83+9=92'),\n  @('31-28=3', '64-53=11'),\n  @('43-43=0', '35+55=90'),\n  @('45+27=72', '37+61=98'),\n  @('90-57=33', '67+17=84'),\n  @('56-12=44', '68-16=52'),\n  @('89-51=38', '66+33=99'),\n  @('36+48=84', '21+41=62'),\n  @('94-63=31', '46+50=96'),\n  @('71-40=31', '29+59=88'),\n  @('6+77=83', '38+61=99'),\n  @('9+54=63', '99-49=50'),\n  @('28-21=7', '37-8=29'),\n  @('75-53=22', '29-11=18'),\n  @('64+9=73', '11+42=53'),\n  @('5+70=75', '59+7=66'),\n  @('41-33=8', '58+10=68'),\n)\n\n$d = $word.ActiveDocument\n\n# Title paragraph (paragraph 1, not part of the table).\n$titlePara = $d.Paragraphs(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r\") -eq $DateOld) {\n  $titlePara.Range.Text = $DateNew\n}\n\n$tbl = $d.Tables(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $pair = $CellPairs[$i]\n    $old = $pair[0]\n    $new = $pair[1]\n    $cell = $tbl.Cell($r, $c)\n    $cellText = $cell.Range.Text.TrimEnd(\"`a\").TrimEnd(\"`r\")\n    if ($cellText -eq $old) {\n      $cell.Range.Text = $new\n    }\n    $i++\n  }\n}\n"}
